$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing "人名" /
# "電話分機號碼" columns from A/B to B/C.
$ws.Columns("A:A").Insert()

# Populate the new column A with the "類別" (Category) header and mark
# the two "RD" rows (Joey / Jack).
$ws.Range("A1").Value = "類別"
$ws.Range("A2").Value = "RD"
$ws.Range("A21").Value = "RD"

# Restore the cell selection recorded in the saved workbook.
$ws.Range("D13").Select() | Out-Null
